# Updates the "Price" (column D) and, for two rows, the "Volume(1h)"
# (column E) figures to match the refreshed coinranking.com scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'24.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.286"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05773"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'6.485"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'3.145"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8170"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8572"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1364"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.06950"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Value = "'0.02896"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09399"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.760"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001526"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04667"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0006005"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006177"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.001237"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.004616"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20HotbitTokenHTBBestin24h"
$ws.Range("D22").Value = "'0.00006106"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Value = "'2.149"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3196"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1285"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.1328"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'0.0002334"
$ws.Range("D28").Style = "Normal"
$ws.Range("D41").Value = "'0.006286"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").Value = "'0.1056"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.002773"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.008519"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005272"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.4404"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.002314"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("D50").Style = "Normal"
